$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.618.43"
$ws.Range("E2").Value = "  +1.23%  "

$ws.Range("D3").Value = "3.025.15"
$ws.Range("E3").Value = "  +2.64%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "379.12"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.93"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.91%  "

$ws.Range("E7").Value = "  +1.17%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  +2.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.56"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.10%  "

$ws.Range("E11").Value = "  -0.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0860"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.27%  "

$ws.Range("D13").Value = "3.510.54"
$ws.Range("E13").Value = "  +3.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.51"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.12%  "

$ws.Range("E15").Value = "  +0.55%  "

$ws.Range("D16").Value = "3.045.94"
$ws.Range("E16").Value = "  +2.94%  "

$ws.Range("E17").Value = "  -2.02%  "

$ws.Range("E18").Value = "  -10.79%  "

$ws.Range("D19").Value = "51.624.53"
$ws.Range("E19").Value = "  +1.29%  "

$ws.Range("E20").Value = "  -0.30%  "

$ws.Range("E21").Value = "  +0.39%  "

$ws.Range("E22").Value = "  +1.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.95"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.66"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.68%  "

$ws.Range("E25").Value = "  -2.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.22"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.31%  "

$ws.Range("E27").Value = "  +6.59%  "

$ws.Range("E28").Value = "  +5.28%  "

$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.07%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.28"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.61%  "

$ws.Range("E31").Value = "  +0.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.27"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.13%  "

$ws.Range("E33").Value = "  +1.71%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "50.59"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.19%  "

$ws.Range("E35").Value = "  +4.58%  "

$ws.Range("E36").Value = "  -0.27%  "

$ws.Range("E38").Value = "  +7.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.289"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +10.36%  "

$ws.Range("E40").Value = "  +3.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.86"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.03%  "

$ws.Range("E42").Value = "  +2.75%  "

$ws.Range("E43").Value = "  -0.53%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.95"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.99%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.72"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +7.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.76"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.40%  "

$ws.Range("E47").Value = "  +3.61%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.39"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.06%  "

$ws.Range("D49").Value = "2.032.18"
$ws.Range("E49").Value = "  +1.10%  "

$ws.Range("D50").Value = "3.323.81"
$ws.Range("E50").Value = "  +2.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0320"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.61%  "
